# RACP update and webapp data text updates
#
# This script applies the meaningful content/view changes captured in the
# commit:
#   - RACP!B2 "Cost Cap" value updated from 180 to 999
#   - The active/selected worksheet tab switches from "About" to "RACP"
#   - The selection on the RACP sheet moves to E8

$wb = $excel.ActiveWorkbook

$wsRACP = $wb.Worksheets.Item("RACP")

# Update the RPS Alternative Compliance Payment cost cap value.
$wsRACP.Range("B2").Value = 999

# Make RACP the active sheet (was "About") and move the selection there.
$wsRACP.Activate() | Out-Null
$wsRACP.Range("E8").Select() | Out-Null
